$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$block1 = New-Object 'object[,]' 24,6
$block1[0,0] = 12.918295161877
$block1[0,1] = 7.103945128901083
$block1[0,2] = 13.106672879016
$block1[0,3] = 35.58975141445249
$block1[0,4] = 48.48959649390691
$block1[0,5] = 19.06744198685161
$block1[1,0] = 12.89303684160258
$block1[1,1] = 7.127969977274086
$block1[1,2] = 13.12951576712277
$block1[1,3] = 35.59459545528607
$block1[1,4] = 48.47125231731943
$block1[1,5] = 19.1240674005624
$block1[2,0] = 12.88053874046285
$block1[2,1] = 7.143363264747961
$block1[2,2] = 13.14589474471633
$block1[2,3] = 35.60909272466211
$block1[2,4] = 48.47803157103592
$block1[2,5] = 19.16330571935821
$block1[3,0] = 12.87620540207645
$block1[3,1] = 7.149798315937336
$block1[3,2] = 13.15316015756776
$block1[3,3] = 35.61788741569617
$block1[3,4] = 48.4853125045975
$block1[3,5] = 19.1804151840398
$block1[4,0] = 12.87553180615337
$block1[4,1] = 7.150876667287077
$block1[4,2] = 13.15440223069518
$block1[4,3] = 35.61952177697255
$block1[4,4] = 48.48679376113039
$block1[4,5] = 19.18332368000112
$block1[5,0] = 12.8804772203395
$block1[5,1] = 7.143449392557521
$block1[5,2] = 13.14599033759708
$block1[5,3] = 35.60919965970772
$block1[5,4] = 48.47811149790135
$block1[5,5] = 19.16353193679377
$block1[6,0] = 12.90896306255167
$block1[6,1] = 7.112096110438987
$block1[6,2] = 13.11406013947048
$block1[6,3] = 35.58902432560549
$block1[6,4] = 48.47951682632047
$block1[6,5] = 19.08603599290167
$block1[7,0] = 12.988563597117
$block1[7,1] = 7.055671935964545
$block1[7,2] = 13.07015645676036
$block1[7,3] = 35.64126264903014
$block1[7,4] = 48.62605229249938
$block1[7,5] = 18.96973113091391
$block1[8,0] = 13.06127795791999
$block1[8,1] = 7.017253896737169
$block1[8,2] = 13.04935597988921
$block1[8,3] = 35.73596464366022
$block1[8,4] = 48.82189571095721
$block1[8,5] = 18.90628449656656
$block1[9,0] = 13.09738249617869
$block1[9,1] = 7.000425831450706
$block1[9,2] = 13.04238913344264
$block1[9,3] = 35.79130196128144
$block1[9,4] = 48.93014881395691
$block1[9,5] = 18.88225113467887
$block1[10,0] = 13.11148291078159
$block1[10,1] = 6.994145969235007
$block1[10,2] = 13.04011027382665
$block1[10,3] = 35.814017052692
$block1[10,4] = 48.97388964074673
$block1[10,5] = 18.87384871784806
$block1[11,0] = 13.10842719540621
$block1[11,1] = 6.995494344184033
$block1[11,2] = 13.04058507887873
$block1[11,3] = 35.80904673954651
$block1[11,4] = 48.96434724441883
$block1[11,5] = 18.87562719946914
$block1[12,0] = 13.09853400094301
$block1[12,1] = 6.99990733228719
$block1[12,2] = 13.04219444627218
$block1[12,3] = 35.79313549488663
$block1[12,4] = 48.93369239270118
$block1[12,5] = 18.88154584104583
$block1[13,0] = 13.09252971294793
$block1[13,1] = 7.002622448670233
$block1[13,2] = 13.04322703891795
$block1[13,3] = 35.78361849089806
$block1[13,4] = 48.91527292566303
$block1[13,5] = 18.88526226335491
$block1[14,0] = 13.05897877531683
$block1[14,1] = 7.018366640436051
$block1[14,2] = 13.04986153882227
$block1[14,3] = 35.73259470733986
$block1[14,4] = 48.81520616541081
$block1[14,5] = 18.90795265811104
$block1[15,0] = 13.03916687672105
$block1[15,1] = 7.028190792439346
$block1[15,2] = 13.05457111319326
$block1[15,3] = 35.70443161458347
$block1[15,4] = 48.75872254082394
$block1[15,5] = 18.9231122513174
$block1[16,0] = 13.02805690472259
$block1[16,1] = 7.033902465031224
$block1[16,2] = 13.05751479190883
$block1[16,3] = 35.68938693968911
$block1[16,4] = 48.7280391573114
$block1[16,5] = 18.93228592714488
$block1[17,0] = 13.02434445714508
$block1[17,1] = 7.03584685098726
$block1[17,2] = 13.05855179051104
$block1[17,3] = 35.68449130333203
$block1[17,4] = 48.71796037302333
$block1[17,5] = 18.93546987021357
$block1[18,0] = 13.0412464059646
$block1[18,1] = 7.027138678209203
$block1[18,2] = 13.05404546028846
$block1[18,3] = 35.70731020026435
$block1[18,4] = 48.76454860134299
$block1[18,5] = 18.92145144035084
$block1[19,0] = 13.10142830230103
$block1[19,1] = 6.998608624108556
$block1[19,2] = 13.04171198045583
$block1[19,3] = 35.79776127843755
$block1[19,4] = 48.94262198110039
$block1[19,5] = 18.87978840377094
$block1[20,0] = 13.14325415682783
$block1[20,1] = 6.980501784207932
$block1[20,2] = 13.03574595765026
$block1[20,3] = 35.86713223451093
$block1[20,4] = 49.07501187019739
$block1[20,5] = 18.85663233275616
$block1[21,0] = 13.12070516033584
$block1[21,1] = 6.990116635315189
$block1[21,2] = 13.03873834476556
$block1[21,3] = 35.82917071564034
$block1[21,4] = 49.00289200929786
$block1[21,5] = 18.86861716482528
$block1[22,0] = 13.04030537766774
$block1[22,1] = 7.02761414075577
$block1[22,2] = 13.05428237264669
$block1[22,3] = 35.70600521974779
$block1[22,4] = 48.76190906376087
$block1[22,5] = 18.92220086594208
$block1[23,0] = 12.96451011574478
$block1[23,1] = 7.070399462594007
$block1[23,2] = 13.08002494473775
$block1[23,3] = 35.6172561153225
$block1[23,4] = 48.57094364934061
$block1[23,5] = 18.99734899921087
$ws.Range("C2:H25").Value = $block1

$block2 = New-Object 'object[,]' 24,3
$block2[0,0] = 16.92297178975247
$block2[0,1] = 9.041396333130493
$block2[0,2] = 18.99571619886904
$block2[1,0] = 16.52711855544097
$block2[1,1] = 9.064217242582453
$block2[1,2] = 18.84518672540883
$block2[2,0] = 16.28253997284509
$block2[2,1] = 9.079128988207758
$block2[2,2] = 18.75591056551995
$block2[3,0] = 16.18263658412292
$block2[3,1] = 9.085432266825419
$block2[3,2] = 18.72035403059067
$block2[4,0] = 16.16603788089413
$block2[4,1] = 9.086492620789555
$block2[4,2] = 18.7145005901843
$block2[5,0] = 16.2811933962777
$block2[5,1] = 9.079213078313456
$block2[5,2] = 18.75542765890399
$block2[6,0] = 16.78688440083184
$block2[6,1] = 9.049078511499676
$block2[6,2] = 18.94317971773953
$block2[7,0] = 17.76028493596618
$block2[7,1] = 8.997104460543927
$block2[7,2] = 19.33487188304712
$block2[8,0] = 18.45653557215127
$block2[8,1] = 8.963234217349928
$block2[8,2] = 19.63492397590771
$block2[9,0] = 18.76768279774286
$block2[9,1] = 8.948757458716704
$block2[9,2] = 19.77363589346464
$block2[10,0] = 18.88459082440583
$block2[10,1] = 8.943408971679773
$block2[10,2] = 19.82644218917697
$block2[11,0] = 18.85945508472485
$block2[11,1] = 8.944554929648918
$block2[11,2] = 19.81505762914014
$block2[12,0] = 18.77731988631687
$block2[12,1] = 8.94831476083497
$block2[12,2] = 19.77797493350619
$block2[13,0] = 18.7268870861661
$block2[13,1] = 8.950635148193621
$block2[13,2] = 19.75529587816633
$block2[14,0] = 18.43607911431458
$block2[14,1] = 8.964199013248544
$block2[14,2] = 19.62590004778721
$block2[15,0] = 18.25616452160415
$block2[15,1] = 8.972758217791453
$block2[15,2] = 19.54705952574392
$block2[16,0] = 18.15216340820511
$block2[16,1] = 8.97776889374823
$block2[16,2] = 19.50192358756158
$block2[17,0] = 18.11686490006816
$block2[17,1] = 8.979480486205869
$block2[17,2] = 19.48667878539856
$block2[18,0] = 18.27537130601393
$block2[18,1] = 8.971838007048847
$block2[18,2] = 19.55543067810025
$block2[19,0] = 18.8014707305319
$block2[19,1] = 8.947206785416824
$block2[19,2] = 19.7888597624971
$block2[20,0] = 19.1399171547673
$block2[20,1] = 8.931887126657561
$block2[20,2] = 19.9430274796365
$block2[21,0] = 18.95981060338925
$block2[21,1] = 8.939992411653201
$block2[21,2] = 19.86061125601066
$block2[22,0] = 18.26668967240287
$block2[22,1] = 8.972253754347436
$block2[22,2] = 19.55164548362991
$block2[23,0] = 17.49974133222961
$block2[23,1] = 9.010405239005328
$block2[23,2] = 19.22660040362172
$ws.Range("K2:M25").Value = $block2
